$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(16,0,1,0,-0.02589133865184522,393.79443359375,395.720277444279,393.7953796386719,20.85106086730957,21.9060318600473,21.42024993896484,53.80360412597656,74.55619062043633,53.80360412597656,1.64277195930481,3.753246358879458,1.64277195930481,400,21.02386570101565,50.73561603332475,1.63925),
    @(17,0,1,0,-0.02589133865184522,393.79443359375,395.4554605582848,393.7952575683594,20.85170555114746,21.95919121404069,21.4454231262207,53.73841857910156,77.25520396924321,53.73841857910156,1.648123621940613,3.753246358879458,1.648123621940613,400,21.03198133034893,50.80485205926722,1.63925),
    @(18,0,1,0,-0.02589133865184522,393.79443359375,394.9270652328814,393.7949829101562,20.85770416259766,21.90874293348235,21.41815185546875,53.74079132080078,80.05135745566002,53.74079132080078,1.651173949241638,3.753246358879458,1.651173949241638,400,21.03492439835391,50.67465438929317,1.6353),
    @(19,0,1,0,-0.02589133865184522,393.79443359375,394.4952244896864,393.7947692871094,20.85620498657227,21.91777397581138,21.4197940826416,53.72517013549805,82.30836977849175,53.72517013549805,1.64588463306427,3.753246358879458,1.64588463306427,400,21.03449933925149,50.97157221308606,1.6353),
    @(20,0,1,0,-0.02595351192412863,393.79443359375,395.4554561242449,393.7952575683594,20.85307693481445,21.95919121404069,21.43817901611328,53.69776916503906,77.25520396924321,53.69776916503906,1.650556802749634,3.753246358879458,1.650556802749634,400,21.03830330421236,50.69720672778764,1.636616666666667),
    @(21,0,1,0,-0.02595351192412863,393.79443359375,395.1721384114802,393.7951049804688,20.84575653076172,21.98307220568238,21.44550704956055,53.71715545654297,79.90490832020177,53.71715545654297,1.653978824615479,3.753246358879458,1.653978824615479,400,21.03700657083391,50.57297699337158,1.637275),
    @(22,0,1,0,-0.02595351192412863,393.79443359375,394.6595917536392,393.7948608398438,20.86370086669922,21.92500826988255,21.42190170288086,53.63594818115234,82.50714421290616,53.63594818115234,1.65619945526123,3.753246358879458,1.65619945526123,400,21.02300839875534,50.40976523532517,1.636616666666667),
    @(23,0,1,0,-0.02595351192412863,393.79443359375,394.2053980332198,393.7946472167969,20.87457084655762,21.89966522430889,21.41254425048828,53.59566497802734,84.69844944499691,53.59566497802734,1.662029266357422,3.753246358879458,1.662029266357422,400,21.016668455537,50.71032798025819,1.63925),
    @(24,0,1,0,-0.02596851188668779,393.79443359375,395.1721339806171,393.7951049804688,20.8687686920166,21.98307220568238,21.4525089263916,53.57375335693359,79.90490832020178,53.57375335693359,1.661708950996399,3.753246358879458,1.661708950996399,400,21.0306594424687,50.50588047367383,1.63925),
    @(25,0,1,0,-0.02596851188668779,393.79443359375,394.8566111674057,393.7949523925781,20.884033203125,21.97170544887955,21.45294189453125,53.57630157470703,82.51400668820797,53.57630157470703,1.662889719009399,3.753246358879458,1.662889719009399,400,21.00521344540757,50.58283207959067,1.63925),
    @(26,0,1,0,-0.02596851188668779,393.79443359375,394.3665252130789,393.7947082519531,20.88093566894531,21.9140822631777,21.42060470581055,53.50544738769531,84.85573318594672,53.50544738769531,1.666810870170593,3.753246358879458,1.666810870170593,400,21.00612105448866,50.48521509443287,1.63925),
    @(27,0,1,0,-0.02596851188668779,393.79443359375,393.9419071928118,393.7944946289062,20.89557456970215,21.89472215812042,21.41690254211426,53.41579818725586,86.75320228667478,53.41579818725586,1.665271997451782,3.753246358879458,1.665271997451782,400,21.01284293241915,50.28880636117213,1.63925),
)

$startRow = 18
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowVals[$j]
    }
}
